$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A7").Value = "Propiedad "
$ws.Range("B7").Value = "Property"

$ws.Range("B7").HorizontalAlignment = -4131  # xlLeft
$ws.Range("B7").IndentLevel = 1
